$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: region names reshuffled (ranking order changed, labels moved between rows)
$ws.Range("A2").Value = "Distrito Federal"
$ws.Range("A4").Value = "Santa Catarina"
$ws.Range("A5").Value = "Goiás"
$ws.Range("A6").Value = "Rio Grande do Sul"
$ws.Range("A7").Value = "Mato Grosso"

# Column C: quarter reference date updated for every data row.
# Force text format first so Excel keeps "01/04/2025" as a literal string
# instead of auto-converting it into a date serial, then restore the
# default (Normal) style so no stray formatting is left behind.
$cDates = $ws.Range("C2:C10")
$cDates.NumberFormat = "@"
$cDates.Value = "01/04/2025"
$cDates.Style = "Normal"

# Column D: refreshed numeric values
$ws.Range("D2").Value = 55.79
$ws.Range("D3").Value = 55.67
$ws.Range("D4").Value = 55.62
$ws.Range("D5").Value = 54.89
$ws.Range("D6").Value = 54.54
$ws.Range("D7").Value = 54.28
$ws.Range("D8").Value = 44.83
$ws.Range("D9").Value = 51.04
$ws.Range("D10").Value = 44.07
